$d = $word.ActiveDocument

$d.Content.Find.Execute("简历 - Alex Wilber", $true, $false, $false, $false, $false, $true, 1, $false, "简历：Alex Wilber", 2)
$d.Content.Find.Execute("Spark 动画：动画设计器（2021 年 1 月 - 演示）", $true, $false, $false, $false, $false, $true, 1, $false, "Spark Animation：动画设计师（2021 年 1 月 - 今）", 2)
$d.Content.Find.Execute("Pixel Studio：动画设计器（2018 年 6 月 - 2020 年 12 月）", $true, $false, $false, $false, $false, $true, 1, $false, "Pixel Studio：动画设计师（2018 年 6 月 - 2020 年 12 月）", 2)
$d.Content.Find.Execute("闪光动画：初级动画设计师（2016 年 9 月 - 2018 年 5 月）", $true, $false, $false, $false, $false, $true, 1, $false, "Flash Animation：初级动画设计师（2016 年 9 月 - 2018 年 5 月）", 2)
$d.Content.Find.Execute("动画艺术大师，预期毕业：2025年12月", $true, $false, $false, $false, $false, $true, 1, $false, "动画艺术硕士，预计 2025 年 12 月毕业", 2)
$d.Content.Find.Execute("纽约：Spark出版社。", $true, $false, $false, $false, $false, $true, 1, $false, "纽约：Spark Press 出版社。", 2)
